# Auto-generated-style edit script: apply the cell value updates described by the commit diff.
# Each worksheet tab (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets its changed
# H:N "profit calculation" columns refreshed to the latest scheduled-run figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H34").Value = 10332.333
$ws.Range("I34").Value = 997
$ws.Range("K34").Value = 997
$ws.Range("M34").Value = -794
$ws.Range("H36").Value = 10332.333
$ws.Range("I36").Value = 997
$ws.Range("K36").Value = 997
$ws.Range("M36").Value = -282
$ws.Range("H38").Value = 2315.1428
$ws.Range("I38").Value = 219.27272
$ws.Range("K38").Value = 657.81816
$ws.Range("M38").Value = -285.81816
$ws.Range("H40").Value = 3552.4
$ws.Range("J40").Value = 4659.1
$ws.Range("L40").Value = 4659.1
$ws.Range("N40").Value = -5009.1
$ws.Range("H43").Value = 5127.2856
$ws.Range("I43").Value = 3249.75
$ws.Range("K43").Value = 3249.75
$ws.Range("M43").Value = -3180.75
$ws.Range("H100").Value = 1713.909
$ws.Range("J100").Value = 2576
$ws.Range("L100").Value = 2576
$ws.Range("N100").Value = -3658
$ws.Range("H106").Value = 1897.2222
$ws.Range("I106").Value = 1810.7142
$ws.Range("K106").Value = 1810.7142
$ws.Range("M106").Value = -1179.7142
$ws.Range("H112").Value = 3694.6
$ws.Range("J112").Value = 3694.6
$ws.Range("L112").Value = 11083.8
$ws.Range("N112").Value = -13299.8
$ws.Range("H118").Value = 995
$ws.Range("I118").Value = 995
$ws.Range("K118").Value = 2985
$ws.Range("M118").Value = -1328
$ws.Range("H132").Value = 3797.8462
$ws.Range("I132").Value = 811.34485
$ws.Range("K132").Value = 2434.03455
$ws.Range("M132").Value = 95.96545000000015
$ws.Range("H138").Value = 2445.1667
$ws.Range("I138").Value = 2322.6
$ws.Range("J138").Value = 2532.7144
$ws.Range("K138").Value = 6967.799999999999
$ws.Range("L138").Value = 7598.1432
$ws.Range("M138").Value = -1827.799999999999
$ws.Range("N138").Value = -17878.1432

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2945.46
$ws.Range("I32").Value = 2945.46
$ws.Range("K32").Value = 2945.46
$ws.Range("M32").Value = -2658.46
$ws.Range("H45").Value = 7149.3
$ws.Range("I45").Value = 3375
$ws.Range("J45").Value = 9665.5
$ws.Range("K45").Value = 3375
$ws.Range("L45").Value = 9665.5
$ws.Range("M45").Value = -2998
$ws.Range("N45").Value = -10419.5
$ws.Range("H74").Value = 4104.3706
$ws.Range("I74").Value = 3536.353
$ws.Range("J74").Value = 5070
$ws.Range("K74").Value = 3536.353
$ws.Range("L74").Value = 5070
$ws.Range("M74").Value = -2662.353
$ws.Range("N74").Value = -6818
$ws.Range("H77").Value = 4104.3706
$ws.Range("I77").Value = 3536.353
$ws.Range("J77").Value = 5070
$ws.Range("K77").Value = 17681.765
$ws.Range("L77").Value = 25350
$ws.Range("M77").Value = -13313.765
$ws.Range("N77").Value = -34086

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()  # was -5994
$ws.Range("H82").Value = 33333.332
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 45000
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 33333.332
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 45000
$ws.Range("N85").Value = -47652
$ws.Range("H94").Value = 1624.5385
$ws.Range("I94").Value = 1444.1052
$ws.Range("J94").Value = 2114.2856
$ws.Range("K94").Value = 1444.1052
$ws.Range("L94").Value = 2114.2856
$ws.Range("M94").Value = -993.1052
$ws.Range("N94").Value = -3016.2856
$ws.Range("H134").Value = 3567.4546
$ws.Range("I134").Value = 3567.4546
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10702.3638
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8167.363799999999
$ws.Range("N134").ClearContents()  # was -8076

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H68").Value = 41544
$ws.Range("J68").Value = 41544
$ws.Range("L68").Value = 41544
$ws.Range("N68").Value = -43042
$ws.Range("H71").Value = 41544
$ws.Range("J71").Value = 41544
$ws.Range("L71").Value = 124632
$ws.Range("N71").Value = -132120
$ws.Range("H74").Value = 40569.832
$ws.Range("J74").Value = 40569.832
$ws.Range("L74").Value = 40569.832
$ws.Range("N74").Value = -42317.832
$ws.Range("H77").Value = 40569.832
$ws.Range("J77").Value = 40569.832
$ws.Range("L77").Value = 121709.496
$ws.Range("N77").Value = -130445.496
$ws.Range("H81").Value = 44000
$ws.Range("J81").Value = 44000
$ws.Range("L81").Value = 44000
$ws.Range("N81").Value = -45996
$ws.Range("H84").Value = 44000
$ws.Range("J84").Value = 44000
$ws.Range("L84").Value = 132000
$ws.Range("N84").Value = -141984
$ws.Range("H93").Value = 10166.333
$ws.Range("I93").Value = 10166.333
$ws.Range("K93").Value = 10166.333
$ws.Range("M93").Value = -8294.333000000001
$ws.Range("H107").Value = 536.16327
$ws.Range("I107").Value = 490
$ws.Range("K107").Value = 490
$ws.Range("M107").Value = 1430
$ws.Range("H110").Value = 79590
$ws.Range("J110").Value = 79590
$ws.Range("L110").Value = 79590
$ws.Range("N110").Value = -87770
$ws.Range("H111").Value = 79990
$ws.Range("J111").Value = 79990
$ws.Range("L111").Value = 79990
$ws.Range("N111").Value = -88170
$ws.Range("H112").Value = 75992.25
$ws.Range("J112").Value = 75992.25
$ws.Range("L112").Value = 75992.25
$ws.Range("N112").Value = -78946.25

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H69").Value = 1000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()  # was -4472
$ws.Range("H72").Value = 1000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()  # was -16662
$ws.Range("H140").Value = 2093.818
$ws.Range("I140").Value = 1629
$ws.Range("K140").Value = 4887
$ws.Range("M140").Value = 293

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 14361.385
$ws.Range("I70").Value = 11783
$ws.Range("K70").Value = 11783
$ws.Range("M70").Value = -11513
$ws.Range("H73").Value = 14361.385
$ws.Range("I73").Value = 11783
$ws.Range("K73").Value = 11783
$ws.Range("M73").Value = -10847
$ws.Range("H99").Value = 16096
$ws.Range("I99").Value = 10620
$ws.Range("J99").Value = 38000
$ws.Range("K99").Value = 10620
$ws.Range("L99").Value = 38000
$ws.Range("M99").Value = -8374
$ws.Range("N99").Value = -42492

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1594.375
$ws.Range("I22").Value = 2293.3333
$ws.Range("J22").Value = 1175
$ws.Range("K22").Value = 2293.3333
$ws.Range("L22").Value = 1175
$ws.Range("M22").Value = -1998.3333
$ws.Range("N22").Value = -1765
$ws.Range("H27").Value = 1594.375
$ws.Range("I27").Value = 2293.3333
$ws.Range("J27").Value = 1175
$ws.Range("K27").Value = 2293.3333
$ws.Range("L27").Value = 1175
$ws.Range("M27").Value = -2186.3333
$ws.Range("N27").Value = -1389
$ws.Range("H46").Value = 11266.617
$ws.Range("I46").Value = 5161.727
$ws.Range("J46").Value = 14186.348
$ws.Range("K46").Value = 5161.727
$ws.Range("L46").Value = 14186.348
$ws.Range("M46").Value = -4973.727
$ws.Range("N46").Value = -14562.348
$ws.Range("H61").Value = 60827.883
$ws.Range("I61").Value = 79022.38
$ws.Range("J61").Value = 1695.75
$ws.Range("K61").Value = 79022.38
$ws.Range("L61").Value = 1695.75
$ws.Range("M61").Value = -78820.38
$ws.Range("N61").Value = -2099.75
$ws.Range("H101").Value = 72244.75
$ws.Range("J101").Value = 72244.75
$ws.Range("L101").Value = 72244.75
$ws.Range("N101").Value = -78734.75
$ws.Range("H113").Value = 60827.883
$ws.Range("I113").Value = 79022.38
$ws.Range("J113").Value = 1695.75
$ws.Range("K113").Value = 79022.38
$ws.Range("L113").Value = 1695.75
$ws.Range("M113").Value = -76852.38
$ws.Range("N113").Value = -6035.75
$ws.Range("H122").Value = 3472.5
$ws.Range("I122").Value = 3430.7407
$ws.Range("K122").Value = 10292.2221
$ws.Range("M122").Value = -7842.222099999999

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H70").Value = 36625
$ws.Range("J70").Value = 40500
$ws.Range("L70").Value = 40500
$ws.Range("N70").Value = -41130
$ws.Range("H73").Value = 36625
$ws.Range("J73").Value = 40500
$ws.Range("L73").Value = 40500
$ws.Range("N73").Value = -42684
$ws.Range("H107").Value = 324.6
$ws.Range("I107").Value = 221.14285
$ws.Range("K107").Value = 663.4285500000001
$ws.Range("M107").Value = 1256.57145
$ws.Range("H113").Value = 690.63635
$ws.Range("I113").Value = 690.63635
$ws.Range("K113").Value = 2071.90905
$ws.Range("M113").Value = 98.09094999999979
$ws.Range("H126").Value = 2222.25
$ws.Range("I126").Value = 2037.0667
$ws.Range("K126").Value = 6111.2001
$ws.Range("M126").Value = -3641.2001
